# This script reproduces the data refresh described by the commit
# 'Update gh-pages to output generated at 456a3b4': it bumps the
# '想去人数' (interest-count) figures that bilibili re-scraped, rolls the
# 展览 (Exhibitions) sheet's four-event window forward by one slot
# (dropping the event that has passed, shifting the remaining three up,
# and appending the newly announced one), and appends a freshly
# announced listing to the 本地生活 (Local Life) sheet.

$wb = $excel.ActiveWorkbook

$wsExpo  = $wb.Worksheets.Item("展览")
$wsShow  = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll   = $wb.Worksheets.Item("全部类型")

# --- 展览: refreshed '想去人数' (interest-count) values ---
$wsExpo.Range("F2").Value = 8986
$wsExpo.Range("F3").Value = 1971
$wsExpo.Range("F4").Value = 6606
$wsExpo.Range("F5").Value = 175
$wsExpo.Range("F6").Value = 2135
$wsExpo.Range("F7").Value = 598
$wsExpo.Range("F8").Value = 78
$wsExpo.Range("F10").Value = 70
$wsExpo.Range("F13").Value = 5
$wsExpo.Range("F14").Value = 82
$wsExpo.Range("F15").Value = 20
$wsExpo.Range("F16").Value = 8889
$wsExpo.Range("F17").Value = 168
$wsExpo.Range("F21").Value = 1839
$wsExpo.Range("F22").Value = 862
$wsExpo.Range("F25").Value = 84
$wsExpo.Range("F28").Value = 1035
$wsExpo.Range("F29").Value = 12
$wsExpo.Range("F30").Value = 69
$wsExpo.Range("F31").Value = 550
$wsExpo.Range("F32").Value = 28
$wsExpo.Range("F34").Value = 540
$wsExpo.Range("F35").Value = 2310
$wsExpo.Range("F36").Value = 874
$wsExpo.Range("F37").Value = 538
$wsExpo.Range("F41").Value = 287
$wsExpo.Range("F42").Value = 180
$wsExpo.Range("F44").Value = 1052
$wsExpo.Range("F49").Value = 3995

# --- 演出: refreshed '想去人数' (interest-count) values ---
$wsShow.Range("F6").Value = 3
$wsShow.Range("F14").Value = 13
$wsShow.Range("F26").Value = 10

# --- 本地生活: refreshed '想去人数' (interest-count) values ---
$wsLocal.Range("F3").Value = 722
$wsLocal.Range("F4").Value = 330

# --- 全部类型: refreshed '想去人数' (interest-count) values ---
$wsAll.Range("F3").Value = 8986
$wsAll.Range("F4").Value = 330
$wsAll.Range("F5").Value = 1971
$wsAll.Range("F6").Value = 6606
$wsAll.Range("F7").Value = 175
$wsAll.Range("F8").Value = 2135
$wsAll.Range("F11").Value = 598
$wsAll.Range("F15").Value = 70
$wsAll.Range("F18").Value = 82
$wsAll.Range("F19").Value = 8889
$wsAll.Range("F20").Value = 168
$wsAll.Range("F23").Value = 1839
$wsAll.Range("F24").Value = 862
$wsAll.Range("F26").Value = 84
$wsAll.Range("F28").Value = 1035
$wsAll.Range("F29").Value = 12
$wsAll.Range("F30").Value = 69
$wsAll.Range("F32").Value = 550
$wsAll.Range("F33").Value = 28
$wsAll.Range("F35").Value = 540
$wsAll.Range("F36").Value = 2310
$wsAll.Range("F37").Value = 874
$wsAll.Range("F38").Value = 13
$wsAll.Range("F40").Value = 538
$wsAll.Range("F41").Value = 287
$wsAll.Range("F42").Value = 180
$wsAll.Range("F43").Value = 79
$wsAll.Range("F44").Value = 3995
$wsAll.Range("F46").Value = 10

# --- 展览: rows 45-48 slide forward by one event -----------------------
# Row 45's old entry (New World动漫博览会) has concluded and drops off; the
# three events that follow it each move up one row, and a brand-new
# event (ET金色齿轮国乙同人 only) is appended as the new row 48.

# Row 45
$wsExpo.Range("C45").Value = '杭州·岚梦国潮·夏日盛会'
$wsExpo.Range("D45").Value = '景兴路896号 EFCLIVE欧美广场'
$wsExpo.Range("E45").Value = '2024.11.09 10:00-11.10 18:00'
$wsExpo.Range("F45").Value = 99
$wsExpo.Range("G45").Value = 60
$wsExpo.Range("H45").Value = 'https://show.bilibili.com/platform/detail.html?id=89829'
$wsExpo.Range("I45").Value = '//i0.hdslb.com/bfs/openplatform/202407/t5Yy5W5F1721806075553.jpeg'

# Row 46
$wsExpo.Range("C46").Value = '杭州·巨人only同人展中学篇'
$wsExpo.Range("D46").Value = '康候圣街99号 顺丰创新中心'
$wsExpo.Range("E46").Value = '2024.11.09 09:30-11.09 17:30'
$wsExpo.Range("F46").Value = 16
$wsExpo.Range("G46").Value = 79
$wsExpo.Range("H46").Value = 'https://show.bilibili.com/platform/detail.html?id=92439'
$wsExpo.Range("I46").Value = '//i2.hdslb.com/bfs/openplatform/202409/otLmkybJ1726115788486.jpeg'

# Row 47
$wsExpo.Range("B47").NumberFormat = "@"
$wsExpo.Range("B47").Value = '2024-11-10'
$wsExpo.Range("B47").ClearFormats()
$wsExpo.Range("C47").Value = '杭州·崩坏同人ONLY 爱莉希雅生日会'
$wsExpo.Range("E47").Value = '2024.11.10 08:00-11.10 20:00'
$wsExpo.Range("F47").Value = 79
$wsExpo.Range("H47").Value = 'https://show.bilibili.com/platform/detail.html?id=92228'
$wsExpo.Range("I47").Value = '//i0.hdslb.com/bfs/openplatform/202409/1FsO31h71725897488610.jpeg'

# Row 48
$wsExpo.Range("B48").NumberFormat = "@"
$wsExpo.Range("B48").Value = '2024-11-16'
$wsExpo.Range("B48").ClearFormats()
$wsExpo.Range("C48").Value = '杭州·ET金色齿轮国乙同人only'
$wsExpo.Range("D48").Value = '转塘街道珊瑚沙东路9号 杭州白金汉爵大酒店'
$wsExpo.Range("E48").Value = '2024.11.16 09:30-11.16 22:00'
$wsExpo.Range("F48").Value = 0
$wsExpo.Range("G48").Value = 25
$wsExpo.Range("H48").Value = 'https://show.bilibili.com/platform/detail.html?id=92511'
$wsExpo.Range("I48").Value = '//i1.hdslb.com/bfs/openplatform/202409/XfT00A611726134427042.jpeg'

# --- 本地生活: append newly announced listing as row 5 --------------------
$wsLocal.Range("B5").NumberFormat = "@"
$wsLocal.Range("B5").Value = '2024-09-28'
$wsLocal.Range("C5").Value = '杭州·GOGOGOODS谷子快跑 GBC谷子限量预售（免费入场）'
$wsLocal.Range("D5").Value = '莫干山路隐秀路交叉口 杭州大悦城'
$wsLocal.Range("E5").NumberFormat = "@"
$wsLocal.Range("E5").Value = '2024.09.28 10:00-11.03 22:00'
$wsLocal.Range("F5").Value = 2
$wsLocal.Range("G5").Value = 35
$wsLocal.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=92504'
$wsLocal.Range("I5").Value = '//i2.hdslb.com/bfs/openplatform/202409/mD0obw7u1726286555313.jpeg'

# Copy row 4's cell formatting (border/alignment on column A, plain
# for the rest) down onto the new row 5, then restore A5's numbering,
# which PasteSpecial would otherwise overwrite with row 4's literal 3.
$wsLocal.Range("A4:I4").Copy()
$wsLocal.Range("A5").PasteSpecial(-4122)
$wsLocal.Range("A5").Value = 4

$excel.CutCopyMode = 0

